$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update A2 (value + drop highlight style -> Normal) and B2 (value only)
$ws.Range("A2").Value = 2310429
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 11185

# Row 3: update A3 and B3 values (styles unchanged)
$ws.Range("A3").Value = 2316494
$ws.Range("B3").Value = 30605

# Rows 4-7: clear all values, keep existing formatting/styles
$ws.Range("A4:E7").ClearContents()

# Remove the duplicate-values conditional formatting previously applied to A2
$ws.Range("A2").FormatConditions.Delete()

# Update the active selection shown when the sheet is opened
$ws.Range("A2").Select()
